# Add Config 6 Binary Relevance results for the Ubuntu sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ubuntu")

$ws.Range("C51").Value = "0.561 0.556 0.270 0.199 0.728"
$ws.Range("D51").Value = "0.645 0.204 0.103 0.058 0.918 "
$ws.Range("E51").Value = "0.402 0.468 0.159 0.112 0.573"
$ws.Range("F51").Value = "0.798 0.652 0.849 0.853 0.980"
$ws.Range("C52").Value = "0.748 0.676 0.333 0.546 0.728"
$ws.Range("D52").Value = "0.549 0.442 0.299 0.341 0.931"
$ws.Range("E52").Value = "0.691 0.549 0.201 0.382 0.573"
$ws.Range("F52").Value = "0.784 0.831 0.914 0.926 0.980 "
$ws.Range("C53").Value = "0.694 0.584 0.403 0.534 0.809"
$ws.Range("D53").Value = "0.621 0.540 0.400 0.431 0.755 "
$ws.Range("E53").Value = "0.571 0.424 0.254 0.368 0.684"
$ws.Range("F53").Value = "0.808 0.862 0.923 0.939 0.977 "
$ws.Range("C54").Value = "0.413 0.232 0.110 0.051 0.714"
$ws.Range("D54").Value = "0.560 0.500 0.550 0.333 0.890"
$ws.Range("E54").Value = "0.265 0.132 0.058 0.026 0.556"
$ws.Range("F54").Value = "0.767 0.853 0.933 0.944 0.978"
$ws.Range("C55").Value = "0.533 0.347 0.191 0.273 0.755"
$ws.Range("D55").Value = "0.818 0.945 1.000 0.857 0.947"
$ws.Range("E55").Value = "0.367 0.210 0.106 0.158 0.607"
$ws.Range("F55").Value = "0.824 0.882 0.939 0.953 0.982"
$ws.Range("C57").Value = "0.561 0.556 0.270 0.199 0.728"
$ws.Range("D57").Value = "0.645 0.204 0.103 0.058 0.918"
$ws.Range("E57").Value = "0.402 0.468 0.159 0.112 0.573 "
$ws.Range("F57").Value = "0.798 0.652 0.849 0.853 0.980"
$ws.Range("C58").Value = "0.748 0.676 0.333 0.546 0.728 "
$ws.Range("D58").Value = "0.549 0.442 0.299 0.341 0.931"
$ws.Range("E58").Value = "0.691 0.549 0.201 0.382 0.573"
$ws.Range("F58").Value = "0.784 0.831 0.914 0.926 0.980"
$ws.Range("C59").Value = "0.694 0.584 0.403 0.534 0.809"
$ws.Range("D59").Value = "0.621 0.540 0.400 0.431 0.755"
$ws.Range("E59").Value = "0.571 0.424 0.254 0.368 0.684"
$ws.Range("F59").Value = "0.808 0.862 0.923 0.939 0.977 "
$ws.Range("C60").Value = "0.413 0.232 0.110 0.051 0.714"
$ws.Range("D60").Value = "0.560 0.500 0.550 0.333 0.890"
$ws.Range("E60").Value = "0.265 0.132 0.058 0.026 0.556 "
$ws.Range("F60").Value = "0.767 0.853 0.933 0.944 0.978"
$ws.Range("C61").Value = "0.533 0.347 0.191 0.273 0.755"
$ws.Range("D61").Value = "0.818 0.945 1.000 0.857 0.947"
$ws.Range("E61").Value = "0.367 0.210 0.106 0.158 0.607"
$ws.Range("F61").Value = "0.824 0.882 0.939 0.953 0.982"

# Move the selection cursor to where the user left off editing
$ws.Activate()
$ws.Range("D63").Select()
